$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows of Leetcode entries.
$ws.Range("A43").Value = 165
$ws.Range("B43").Value = "比较版本号"
$ws.Range("D43").Value = "||"
$ws.Range("E43").Value = "思路简单"

$ws.Range("A44").Value = 166
$ws.Range("B44").Value = "分数到小数"
$ws.Range("D44").Value = "||||"
$ws.Range("E44").Value = "不擅长这种题目，要看一看"

# Scroll the view and move the selection to match the edited sheet state.
$excel.ActiveWindow.ScrollRow = 14
$ws.Range("E44").Select()
